# Updated cryptos list on Sun Dec  3 23:07:58 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Price, Volume(1h)) updates. $null means "leave Price unchanged".
$updates = @{
    2  = @('40.129.83', '  +1.88%  ')
    3  = @('2.206.16', '  +2.11%  ')
    4  = @($null, '  +0.08%  ')
    5  = @('228.64', '  -0.12%  ')
    6  = @($null, '  +1.43%  ')
    7  = @('64.03', '  +1.53%  ')
    8  = @($null, '  +0.06%  ')
    9  = @('0.399', '  +0.89%  ')
    10 = @($null, '  -0.09%  ')
    11 = @($null, '  +0.42%  ')
    12 = @('16.03', '  +0.04%  ')
    13 = @('2.533.21', '  +2.12%  ')
    14 = @('22.21', '  +0.35%  ')
    15 = @($null, '  +0.73%  ')
    16 = @('5.60', '  +0.54%  ')
    17 = @('2.204.46', '  +2.07%  ')
    18 = @('40.061.73', '  +1.77%  ')
    19 = @('0.0₃0911', '  +6.89%  ')
    20 = @('72.51', '  +0.45%  ')
    21 = @('6.10', '  -0.38%  ')
    22 = @('232.92', '  +1.97%  ')
    23 = @($null, '  -0.03%  ')
    24 = @('2.33', '  -0.47%  ')
    25 = @('2.36', '  -0.35%  ')
    26 = @('9.73', '  -0.33%  ')
    27 = @('171.87', '  -0.05%  ')
    28 = @('0.141', '  +2.16%  ')
    29 = @($null, '  +2.91%  ')
    30 = @('20.18', '  +2.32%  ')
    31 = @('2.75', '  +6.14%  ')
    32 = @($null, '  +1.41%  ')
    33 = @($null, '  -1.16%  ')
    34 = @('4.75', '  -0.99%  ')
    35 = @('7.08', '  -0.14%  ')
    36 = @('0.0626', '  +0.80%  ')
    37 = @('3.91', '  +9.88%  ')
    38 = @('2.45', '  +0.90%  ')
    41 = @('103.87', '  -0.44%  ')
    42 = @($null, '  -0.48%  ')
    43 = @('17.59', '  -2.51%  ')
    44 = @($null, '  +3.38%  ')
    45 = @('1.522.10', '  -1.10%  ')
    46 = @('8.29', '  +6.77%  ')
    47 = @($null, '  +0.76%  ')
    48 = @('0.0928', '  -0.23%  ')
    49 = @($null, '  -0.46%  ')
    50 = @($null, '  +33.29%  ')
    51 = @('2.410.17', '  +1.96%  ')
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]
    if ($null -ne $price) {
        # Price column holds text that merely looks numeric (e.g. "228.64",
        # thousand-separated "40.129.83", subscript-notation strings, ...).
        # Force text storage so COM doesn't silently coerce it to a number,
        # then restore the default "Normal" style so no new number format
        # is left behind on the cell.
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $price
        $cell.Style = "Normal"
    }
    $ws.Cells.Item($row, 5).Value = $volume
}

# Rows 39 and 40 swap places (FTXToken <-> BinanceUSD) with refreshed values.
$ws.Cells.Item(39, 2).Value = 'FTXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$cell39 = $ws.Cells.Item(39, 4)
$cell39.NumberFormat = "@"
$cell39.Value = '5.04'
$cell39.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +19.87%  '

$ws.Cells.Item(40, 2).Value = 'BinanceUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell40 = $ws.Cells.Item(40, 4)
$cell40.NumberFormat = "@"
$cell40.Value = '1.00'
$cell40.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.12%  '
